$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9999999903147524
$ws.Range("A2").Value = 0.99393559460034719
$ws.Range("A3").Value = 0.97100655283155557
$ws.Range("A4").Value = 0.96434173588988015
$ws.Range("A5").Value = 0.95468548634932704
$ws.Range("A6").Value = 0.93154971137159492
$ws.Range("A7").Value = 0.92862490234993567
$ws.Range("A8").Value = 0.92535009308060212
$ws.Range("A9").Value = 0.92352238969603573
$ws.Range("A10").Value = 0.9226397797478405
$ws.Range("A11").Value = 0.92263451209074354
$ws.Range("A12").Value = 0.91985463138175905
$ws.Range("A13").Value = 0.90856733508727894
$ws.Range("A14").Value = 0.90439999432446783
$ws.Range("A15").Value = 0.90180852627965491
$ws.Range("A16").Value = 0.89930203313715995
$ws.Range("A17").Value = 0.8955941253297397
$ws.Range("A18").Value = 0.89448522743400749
$ws.Range("A19").Value = 0.99053880939568972
$ws.Range("A20").Value = 0.98342183263905159
$ws.Range("A21").Value = 0.98202336471122476
$ws.Range("A22").Value = 0.98075886127756917
$ws.Range("A23").Value = 0.95797823363182233
$ws.Range("A24").Value = 0.94495600208926389
$ws.Range("A25").Value = 0.93849879447323215
$ws.Range("A26").Value = 0.9293258677782037
$ws.Range("A27").Value = 0.92647528177515315
$ws.Range("A28").Value = 0.9141153428918789
$ws.Range("A29").Value = 0.90561817739100148
$ws.Range("A30").Value = 0.9024812011756278
$ws.Range("A31").Value = 0.89482795529184544
$ws.Range("A32").Value = 0.89314869487238691
$ws.Range("A33").Value = 0.8926287037802727
